# cs401_3_2_20.pptx - "Add files via upload" edit
#
# 1. Every cached 'datetimeFigureOut' footer field (slide master, all
#    slide layouts, notes master) advances one day: 3/1/2020 -> 3/2/2020.
# 2. Slide 1: the bullet textbox grows taller to fit a new sub-bullet.
# 3. Slide 1: "assignment" gets curly quotes in the heading bullet, the
#    trailing space is trimmed off the Assignment #2 bullet, and a new
#    sub-bullet explaining the quotation marks is appended.

$p = $ppt.ActivePresentation

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "3/1/2020") {
                $tr.Text = "3/2/2020"
            }
        }
    }
}

# --- 1a. Slide master date placeholder -------------------------------
Update-DateShapes $p.SlideMaster.Shapes

# --- 1b. Every slide layout's date placeholder ------------------------
$sm = $p.SlideMaster
for ($L = 1; $L -le $sm.CustomLayouts.Count; $L++) {
    Update-DateShapes $sm.CustomLayouts.Item($L).Shapes
}

# --- 1c. Notes master date placeholder --------------------------------
# Direct TextRange writes on notes-master shapes don't stick in this
# host, so go through the HeadersFooters.DateAndTime accessor instead.
$nm = $p.NotesMaster
$nmDateText = $null
for ($i = 1; $i -le $nm.Shapes.Count; $i++) {
    $shp = $nm.Shapes.Item($i)
    if ($shp.Name -like "Date Placeholder*") {
        $nmDateText = $shp.TextFrame.TextRange.Text
    }
}
if ($nmDateText -eq "3/1/2020") {
    $nm.HeadersFooters.DateAndTime.Text = "3/2/2020"
}

# --- 2. Grow the bullet textbox on slide 1 to make room for the new ---
#        sub-bullet.
$s1 = $p.Slides.Item(1)
$box = $s1.Shapes.Item(2)
$box.Height = 269.0015748031496

# --- 3. Rewrite / add the bullet text ----------------------------------
$tr = $box.TextFrame.TextRange

$para6 = $tr.Paragraphs(6)
$para6.Runs(1).Text = "Description of the lab “assignment” #7"

$para7 = $tr.Paragraphs(7)
$para7.Runs(1).Text = "This consists of the first part of Assignment #2"

$newBullet = "“assignment” is in quotes because it will not be graded separately as a lab – no demo required"
$tr.InsertAfter("`r(" + $newBullet + ") ") | Out-Null
